# Applies the "Updated symbol list on Sat Dec 31 19:25:18 UTC 2022 with GitHub
# Actions" refresh: refreshed Price (D) / Volume(1h) (E) figures throughout, and
# for the block of rows where a coin fell out of the top-ranked list (9-13 and
# 19-24) each remaining coin/link (B/C) shifted up one row into the vacated slot,
# with the newly-ranked coin written into the row left open at the bottom of the
# block (its own Price/Volume came along with it).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2='246.59', E2='0.57%'
$ws.Range("D2").Value = "'246.59"
$ws.Range("E2").Value = "'0.57%"

# Row 3: D3='26.45', E3='5.17%'
$ws.Range("D3").Value = "'26.45"
$ws.Range("E3").Value = "'5.17%"

# Row 4: D4='5.090', E4='0.39%'
$ws.Range("D4").Value = "'5.090"
$ws.Range("E4").Value = "'0.39%"

# Row 5: D5='0.05599', E5='-0.32%'
$ws.Range("D5").Value = "'0.05599"
$ws.Range("E5").Value = "'-0.32%"

# Row 6: D6='6.484', E6='-0.83%'
$ws.Range("D6").Value = "'6.484"
$ws.Range("E6").Value = "'-0.83%"

# Row 7: D7='0.8129', E7='0.50%'
$ws.Range("D7").Value = "'0.8129"
$ws.Range("E7").Value = "'0.50%"

# Row 8: D8='0.8457', E8='0.78%'
$ws.Range("D8").Value = "'0.8457"
$ws.Range("E8").Value = "'0.78%"

# Row 9: B9='MandalaExchangeToken', C9='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx', D9='0.07003', E9='1.39%'
$ws.Range("B9").Value = "MandalaExchangeToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D9").Value = "'0.07003"
$ws.Range("E9").Value = "'1.39%"

# Row 10: B10='BitrueCoin', C10='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr', D10='0.02852', E10='1.59%'
$ws.Range("B10").Value = "BitrueCoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D10").Value = "'0.02852"
$ws.Range("E10").Value = "'1.59%"

# Row 11: B11='BitMartToken', C11='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx', D11='0.09380', E11='-0.32%'
$ws.Range("B11").Value = "BitMartToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D11").Value = "'0.09380"
$ws.Range("E11").Value = "'-0.32%"

# Row 12: B12='BitForexToken', C12='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf', D12='0.001518', E12='0.66%'
$ws.Range("B12").Value = "BitForexToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D12").Value = "'0.001518"
$ws.Range("E12").Value = "'0.66%"

# Row 13: B13='One', C13='https://coinranking.com/coin/6Lga5NiXX3rT+one-one', D13='0.0006000', E13='0.72%'
$ws.Range("B13").Value = "One"
$ws.Range("C13").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D13").Value = "'0.0006000"
$ws.Range("E13").Value = "'0.72%"

# Row 14: D14='0.006126', E14='-0.01%'
$ws.Range("D14").Value = "'0.006126"
$ws.Range("E14").Value = "'-0.01%"

# Row 15: D15='3.604', E15='2.95%'
$ws.Range("D15").Value = "'3.604"
$ws.Range("E15").Value = "'2.95%"

# Row 16: D16='3.014', E16='0.34%'
$ws.Range("D16").Value = "'3.014"
$ws.Range("E16").Value = "'0.34%"

# Row 17: E17='-1.71%'
$ws.Range("E17").Value = "'-1.71%"

# Row 18: E18='0.76%'
$ws.Range("E18").Value = "'0.76%"

# Row 19: B19='WazirX', C19='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx', D19='0.1338', E19='0.38%'
$ws.Range("B19").Value = "WazirX"
$ws.Range("C19").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D19").Value = "'0.1338"
$ws.Range("E19").Value = "'0.38%"

# Row 20: B20='LiechtensteinCryptoassetsExchange', C20='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx', D20='0.03180', E20='-1.72%'
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "'0.03180"
$ws.Range("E20").Value = "'-1.72%"

# Row 21: B21='ProBitToken', C21='https://coinranking.com/coin/lQP4d6T2+probittoken-prob', D21='0.1297', E21='0.44%'
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1297"
$ws.Range("E21").Value = "'0.44%"

# Row 22: B22='MCDex', C22='https://coinranking.com/coin/3nMM61qeg+mcdex-mcb', D22='3.743', E22='-0.11%'
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "'3.743"
$ws.Range("E22").Value = "'-0.11%"

# Row 23: B23='CoinExToken', C23='https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet', D23='0.04648', E23='-0.84%'
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04648"
$ws.Range("E23").Value = "'-0.84%"

# Row 24: B24='ZBToken', C24='https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb', D24='0.1350', E24='-1.42%'
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").Value = "'0.1350"
$ws.Range("E24").Value = "'-1.42%"

# Row 25: D25='0.001249', E25='0.46%'
$ws.Range("D25").Value = "'0.001249"
$ws.Range("E25").Value = "'0.46%"

# Row 26: D26='0.004582', E26='1.23%'
$ws.Range("D26").Value = "'0.004582"
$ws.Range("E26").Value = "'1.23%"

# Row 27: D27='0.00009600', E27='-0.98%'
$ws.Range("D27").Value = "'0.00009600"
$ws.Range("E27").Value = "'-0.98%"

# Row 28: E28='1.71%'
$ws.Range("E28").Value = "'1.71%"

# Row 40: D40='0.03668', E40='0.92%'
$ws.Range("D40").Value = "'0.03668"
$ws.Range("E40").Value = "'0.92%"

# Row 41: D41='0.006210', E41='84.37%'
$ws.Range("D41").Value = "'0.006210"
$ws.Range("E41").Value = "'84.37%"

# Row 42: D42='0.1056', E42='-22.74%'
$ws.Range("D42").Value = "'0.1056"
$ws.Range("E42").Value = "'-22.74%"

# Row 43: E43='-8.15%'
$ws.Range("E43").Value = "'-8.15%"

# Row 44: D44='0.009245', E44='14.75%'
$ws.Range("D44").Value = "'0.009245"
$ws.Range("E44").Value = "'14.75%"

# Row 45: D45='0.00005276', E45='-0.03%'
$ws.Range("D45").Value = "'0.00005276"
$ws.Range("E45").Value = "'-0.03%"

# Row 46: E46='0.04%'
$ws.Range("E46").Value = "'0.04%"

# Row 47: E47='-42.08%'
$ws.Range("E47").Value = "'-42.08%"

# Row 48: D48='0.002637', E48='28.92%'
$ws.Range("D48").Value = "'0.002637"
$ws.Range("E48").Value = "'28.92%"

# Row 49: E49='0.04%'
$ws.Range("E49").Value = "'0.04%"

# Row 50: E50='0.04%'
$ws.Range("E50").Value = "'0.04%"

